$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.009892702102661
$ws.Range("B1").Value = 2.175233125686646
$ws.Range("C1").Value = 2.231827735900879
$ws.Range("D1").Value = 2.71663236618042
$ws.Range("E1").Value = 3.569836139678955
